$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last 3 rows of data (Trial 17, 18, 19 -> previously rows 18:20)
$ws.Rows("18:20").Delete()

# Add new header + shared string "ITI"
$ws.Range("D1").Value = "ITI"

# Update ConditionType (column C) values for rows 2-17
$cVals = @(1,1,2,2,4,2,1,4,4,3,2,3,1,3,3,4)
for ($i = 0; $i -lt $cVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
}

# Fill new ITI column (column D) values for rows 2-17
$dVals = @(8,6,6,6,6,8,6,9,6,7,6,6,6,9,7,10)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
}

# Column width adjustments (target stored widths: 11.83203125 / 12.33203125 / 14.6640625)
$ws.Columns.Item(3).ColumnWidth = 10.998697916666666
$ws.Columns.Item(5).ColumnWidth = 11.498697916666666
$ws.Columns.Item(6).ColumnWidth = 13.830729166666666

# Update selection to match target view
$ws.Range("J18").Select()
